$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.710.02'
$ws.Range('E2').Value = '  +4.21%  '
$ws.Range('D3').Value = '3.258.17'
$ws.Range('E3').Value = '  +7.43%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '583.11'
$ws.Range('E5').Value = '  +4.97%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '154.46'
$ws.Range('E6').Value = '  +9.43%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '3.250.03'
$ws.Range('E8').Value = '  +7.56%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.515'
$ws.Range('E9').Value = '  +5.45%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.09'
$ws.Range('E10').Value = '  +8.78%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.165'
$ws.Range('E11').Value = '  +6.22%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.490'
$ws.Range('E12').Value = '  +4.96%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '37.98'
$ws.Range('E13').Value = '  +3.53%  '
$ws.Range('E14').Value = '  +5.67%  '
$ws.Range('D15').Value = '3.787.06'
$ws.Range('E15').Value = '  +7.49%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '559.40'
$ws.Range('E16').Value = '  +13.21%  '
$ws.Range('D17').Value = '66.778.47'
$ws.Range('E17').Value = '  +4.07%  '
$ws.Range('D18').Value = '3.256.02'
$ws.Range('E18').Value = '  +7.04%  '
$ws.Range('E19').Value = '  +3.50%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.15'
$ws.Range('E20').Value = '  +6.34%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.50'
$ws.Range('E21').Value = '  +5.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.746'
$ws.Range('E22').Value = '  +8.27%  '
$ws.Range('E23').Value = '  +8.93%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.70'
$ws.Range('E24').Value = '  +7.80%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '82.05'
$ws.Range('E25').Value = '  +3.22%  '
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.28'
$ws.Range('E27').Value = '  +18.23%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.99'
$ws.Range('E28').Value = '  +8.26%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.24'
$ws.Range('E29').Value = '  +6.10%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '27.94'
$ws.Range('E30').Value = '  +7.05%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.77'
$ws.Range('E31').Value = '  +5.11%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.00'
$ws.Range('E32').Value = '  -0.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.18'
$ws.Range('E33').Value = '  +6.49%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '562.17'
$ws.Range('E34').Value = '  +8.93%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.72'
$ws.Range('E35').Value = '  +4.06%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.40'
$ws.Range('E36').Value = '  +7.41%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0459'
$ws.Range('E37').Value = '  +13.14%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '55.38'
$ws.Range('E38').Value = '  +4.85%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0867'
$ws.Range('E39').Value = '  +7.87%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.132'
$ws.Range('E40').Value = '  +7.57%  '
$ws.Range('E41').Value = '  +13.79%  '
$ws.Range('D42').Value = '3.185.93'
$ws.Range('E42').Value = '  +9.77%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.65'
$ws.Range('E43').Value = '  +2.91%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.280'
$ws.Range('E44').Value = '  +13.47%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.33'
$ws.Range('E45').Value = '  +10.34%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '26.46'
$ws.Range('E46').Value = '  +4.55%  '
$ws.Range('B47').Value = 'PEPE'
$ws.Range('C47').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D47').Value = '0.0₃0560'
$ws.Range('E47').Value = '  +3.80%  '
$ws.Range('B48').Value = 'USDe'
$ws.Range('C48').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.999'
$ws.Range('E48').Value = '  +0.08%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '126.35'
$ws.Range('E49').Value = '  +5.22%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.113'
$ws.Range('E50').Value = '  +2.98%  '
$ws.Range('E51').Value = '  +8.63%  '
